$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card20")

# The existing last data row (24) recorded the "766 t / تم تغيير زيت الجيربوكس"
# service event but left columns B:K blank. A new event ("786 t / تم تغيير
# جريده 1 ومعايرها") needs to be appended below it as row 25, and the now
# second-to-last row (the old row 24) gets its blank B:K cells explicitly
# filled with "nan" to match the rest of the table.
#
# Insert a fresh blank row at 24: this pushes the current row 24 (with all
# its existing cells/types/styles intact, including the empty B:K cells)
# down to row 25 - exactly the shape the new row needs before we overwrite
# its L:O columns with the new event.
$ws.Rows.Item(24).Insert(-4121)

# Row 24 is now empty. Re-create its "card" id (A24) by copying the text
# value from A25 so it keeps the original Text type (not re-parsed as a
# number) and no new style gets introduced.
$ws.Cells.Item(25, 1).Copy()
$ws.Cells.Item(24, 1).PasteSpecial(-4163)

# Fill B24:K24 with "nan" (matching every other row in the table).
$ws.Cells.Item(24, 2).Value2 = "nan"
$ws.Cells.Item(24, 3).Value2 = "nan"
$ws.Cells.Item(24, 4).Value2 = "nan"
$ws.Cells.Item(24, 5).Value2 = "nan"
$ws.Cells.Item(24, 6).Value2 = "nan"
$ws.Cells.Item(24, 7).Value2 = "nan"
$ws.Cells.Item(24, 8).Value2 = "nan"
$ws.Cells.Item(24, 9).Value2 = "nan"
$ws.Cells.Item(24, 10).Value2 = "nan"
$ws.Cells.Item(24, 11).Value2 = "nan"

# Restore row 24's Date/Event/Correction/Serviced-by columns (L:O), which
# already held this data before the insert pushed them to row 25.
$ws.Cells.Item(24, 12).Value2 = "14\8\2025"
$ws.Cells.Item(24, 13).Value2 = "766 t"
$ws.Cells.Item(24, 14).Value2 = "تم تغيير زيت الجيربوكس"
$ws.Cells.Item(24, 15).Value2 = "تيم العمل"

# Populate the new row 25 with the new service event. A25:K25 already carry
# the correct card id / blank values from the row-24 shift above.
$ws.Cells.Item(25, 12).Value2 = "27\8\2025"
$ws.Cells.Item(25, 13).Value2 = "786 t"
$ws.Cells.Item(25, 14).Value2 = "تم تغيير جريده 1 ومعايرها"
$ws.Cells.Item(25, 15).Value2 = "الخبير"
